$d = $word.ActiveDocument

$pairs = @(
    @("476÷2=238, 0", "598÷3=199, 1"),
    @("309÷7=44, 1", "612÷4=153, 0"),
    @("359÷3=119, 2", "306÷3=102, 0"),
    @("844÷4=211, 0", "770÷8=96, 2"),
    @("699÷5=139, 4", "761÷4=190, 1"),
    @("894÷4=223, 2", "316÷9=35, 1"),
    @("241÷4=60, 1", "118÷7=16, 6"),
    @("124÷4=31, 0", "536÷7=76, 4"),
    @("181÷7=25, 6", "886÷2=443, 0"),
    @("397÷2=198, 1", "452÷3=150, 2"),
    @("234÷4=58, 2", "920÷3=306, 2"),
    @("196÷4=49, 0", "179÷9=19, 8"),
    @("100÷6=16, 4", "215÷8=26, 7"),
    @("302÷2=151, 0", "162÷6=27, 0"),
    @("791÷2=395, 1", "764÷8=95, 4"),
    @("103÷9=11, 4", "247÷4=61, 3"),
    @("802÷9=89, 1", "317÷3=105, 2"),
    @("480÷7=68, 4", "579÷5=115, 4"),
    @("629÷8=78, 5", "516÷6=86, 0"),
    @("957÷2=478, 1", "724÷9=80, 4"),
    @("678÷3=226, 0", "496÷7=70, 6"),
    @("615÷7=87, 6", "621÷5=124, 1"),
    @("328÷5=65, 3", "914÷7=130, 4"),
    @("272÷5=54, 2", "823÷2=411, 1"),
    @("390÷2=195, 0", "210÷3=70, 0")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
